# Update accelerometer data: replace rows 2-21 with shifted window and
# extend sheet through row 31 with new samples (per "may 9th" edit).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 18.70621109008789
$ws.Cells.Item(2, 2).Value = -5.636280059814453
$ws.Cells.Item(2, 3).Value = 6.914060115814209
$ws.Cells.Item(3, 1).Value = 19.93985939025879
$ws.Cells.Item(3, 2).Value = -3.599599838256836
$ws.Cells.Item(3, 3).Value = 8.073759078979492
$ws.Cells.Item(4, 1).Value = 27.86785316467285
$ws.Cells.Item(4, 2).Value = -18.30047225952148
$ws.Cells.Item(4, 3).Value = 6.377731800079346
$ws.Cells.Item(5, 1).Value = -0.164954662322998
$ws.Cells.Item(5, 2).Value = -10.53660011291504
$ws.Cells.Item(5, 3).Value = 4.150550365447998
$ws.Cells.Item(6, 1).Value = 5.472853660583496
$ws.Cells.Item(6, 2).Value = -46.38559722900391
$ws.Cells.Item(6, 3).Value = -0.9001345634460449
$ws.Cells.Item(7, 1).Value = 3.861607551574707
$ws.Cells.Item(7, 2).Value = -19.79164886474609
$ws.Cells.Item(7, 3).Value = 0.622889518737793
$ws.Cells.Item(8, 1).Value = 14.61344242095947
$ws.Cells.Item(8, 2).Value = -10.30790042877197
$ws.Cells.Item(8, 3).Value = 17.13724899291992
$ws.Cells.Item(9, 1).Value = 2.93750786781311
$ws.Cells.Item(9, 2).Value = -46.01893997192383
$ws.Cells.Item(9, 3).Value = 7.825448036193848
$ws.Cells.Item(10, 1).Value = -20.74547958374023
$ws.Cells.Item(10, 2).Value = -0.2439025640487671
$ws.Cells.Item(10, 3).Value = 8.103152275085449
$ws.Cells.Item(11, 1).Value = -20.94050025939941
$ws.Cells.Item(11, 2).Value = -7.066665649414063
$ws.Cells.Item(11, 3).Value = 6.037558555603027
$ws.Cells.Item(12, 1).Value = 29.87363624572754
$ws.Cells.Item(12, 2).Value = -61.11213684082031
$ws.Cells.Item(12, 3).Value = 14.88786697387695
$ws.Cells.Item(13, 1).Value = 4.961847305297852
$ws.Cells.Item(13, 2).Value = -15.47994041442871
$ws.Cells.Item(13, 3).Value = 18.70905494689941
$ws.Cells.Item(14, 1).Value = -0.8285019397735596
$ws.Cells.Item(14, 2).Value = -6.118541240692139
$ws.Cells.Item(14, 3).Value = 8.952471733093262
$ws.Cells.Item(15, 1).Value = 0.8424484729766846
$ws.Cells.Item(15, 2).Value = -38.40699005126953
$ws.Cells.Item(15, 3).Value = 14.60053634643555
$ws.Cells.Item(16, 1).Value = 1.708237409591675
$ws.Cells.Item(16, 2).Value = 32.82785034179688
$ws.Cells.Item(16, 3).Value = 13.28276348114014
$ws.Cells.Item(17, 1).Value = -19.85161781311035
$ws.Cells.Item(17, 2).Value = -6.960978984832764
$ws.Cells.Item(17, 3).Value = 4.844282150268555
$ws.Cells.Item(18, 1).Value = -47.64518356323242
$ws.Cells.Item(18, 2).Value = -56.76200866699219
$ws.Cells.Item(18, 3).Value = 40.71841812133789
$ws.Cells.Item(19, 1).Value = -2.008986234664917
$ws.Cells.Item(19, 2).Value = -3.396074771881104
$ws.Cells.Item(19, 3).Value = 5.574520111083984
$ws.Cells.Item(20, 1).Value = 1.793292045593261
$ws.Cells.Item(20, 2).Value = 1.589181900024414
$ws.Cells.Item(20, 3).Value = 13.23852920532227
$ws.Cells.Item(21, 1).Value = 0.5213950872421265
$ws.Cells.Item(21, 2).Value = -13.69121932983398
$ws.Cells.Item(21, 3).Value = 13.58244514465332
$ws.Cells.Item(22, 1).Value = -25.09261322021484
$ws.Cells.Item(22, 2).Value = 15.60748481750488
$ws.Cells.Item(22, 3).Value = 0.5675735473632812
$ws.Cells.Item(23, 1).Value = -30.59898567199707
$ws.Cells.Item(23, 2).Value = -12.55906105041504
$ws.Cells.Item(23, 3).Value = 3.974555969238281
$ws.Cells.Item(24, 1).Value = -38.86380767822266
$ws.Cells.Item(24, 2).Value = -84.71040344238281
$ws.Cells.Item(24, 3).Value = 66.18233489990234
$ws.Cells.Item(25, 1).Value = -8.86505126953125
$ws.Cells.Item(25, 2).Value = -2.532943964004517
$ws.Cells.Item(25, 3).Value = 5.462150573730469
$ws.Cells.Item(26, 1).Value = 7.391067981719971
$ws.Cells.Item(26, 2).Value = -1.471791982650757
$ws.Cells.Item(26, 3).Value = 18.45427322387696
$ws.Cells.Item(27, 1).Value = -4.55918025970459
$ws.Cells.Item(27, 2).Value = -21.72416114807129
$ws.Cells.Item(27, 3).Value = -0.08992767333984369
$ws.Cells.Item(28, 1).Value = -17.91468048095703
$ws.Cells.Item(28, 2).Value = 18.05105400085449
$ws.Cells.Item(28, 3).Value = -10.55736446380615
$ws.Cells.Item(29, 1).Value = -62.19514083862305
$ws.Cells.Item(29, 2).Value = -14.345703125
$ws.Cells.Item(29, 3).Value = 3.762966632843018
$ws.Cells.Item(30, 1).Value = -47.85998153686523
$ws.Cells.Item(30, 2).Value = -73.05361175537109
$ws.Cells.Item(30, 3).Value = 52.46365737915039
$ws.Cells.Item(31, 1).Value = -2.905624389648437
$ws.Cells.Item(31, 2).Value = 1.127065658569336
$ws.Cells.Item(31, 3).Value = 7.177680492401123
